$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")
$ws.Activate()

# Update the "Current Version" value from 99.0.4844.51 to 102.0.5005.115
$ws.Range("E2").Value = "102.0.5005.115"

# Move the active selection to F2 (matches the sheetView selection change)
$ws.Range("F2").Select()
